$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("producto") rows 2-29 change from "DESCONOCIDO" to "FRAMBUESA"
$ws.Range("B2:B29").Value = "FRAMBUESA"
